$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (header "Förändrad") holds a date serial that was bumped by one
# day (45180 -> 45181) for every data row (rows 2 through 89).
$ws.Range("C2:C89").Value = 45181
